$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cell H1 ("Save"), matching style of existing header row (B1:G1)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Add new data values in column H for rows 2 and 3
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 0
